# lab3report.docx proofreading pass (grammar/typography fixes), per the
# OOXML diff for this commit. The diff is almost entirely Word's own
# spell/grammar-check markup (<w:proofErr>, run re-splitting, OLE/bookmark
# renumbering) which carries no visible content change; the only real
# textual edits are:
#  - "матриц хранящихся" -> "матриц, хранящихся" (missing comma)
#  - "дает"   -> "даёт"   (ё)
#  - "Причем" -> "Причём" (ё)
#  - "при увеличение" -> "при увеличении" (grammar fix)
#  - "растет" -> "растёт" (ё), both occurrences
#
# These are plain text substitutions; Word re-flows/re-splits the
# surrounding runs on its own, so a scoped Find/Replace over the whole
# document body is sufficient and faithfully reproduces the edit.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "матриц хранящихся в профильном формате" "матриц, хранящихся в профильном формате"

Replace-Text "произведение дает исходную" "произведение даёт исходную"

Replace-Text "Причем y" "Причём y"

Replace-Text "сильно возрастает при увеличение размерности" "сильно возрастает при увеличении размерности"

Replace-Text "с ростом n растет количество" "с ростом n растёт количество"

Replace-Text "суммарная погрешность стремительно растет." "суммарная погрешность стремительно растёт."
